$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 企业数 / v4_cy_scjy_qys  ->  出口贸易额 / v4_cy_my_ck
$ws.Range("A2").Value = "出口贸易额"
$ws.Range("B2").Value = "出口贸易额"
$ws.Range("D2").Value = "v4_cy_my_ck"

# Row 3: 主营业务收入 / v4_cy_scjy_zyyw  ->  进口贸易额 / v4_cy_my_jk
$ws.Range("A3").Value = "进口贸易额"
$ws.Range("B3").Value = "进口贸易额"
$ws.Range("D3").Value = "v4_cy_my_jk"

# Row 4: 利润总额 / v4_cy_scjy_lrze  ->  贸易总额 / v4_cy_my_jck
$ws.Range("A4").Value = "贸易总额"
$ws.Range("B4").Value = "贸易总额"
$ws.Range("D4").Value = "v4_cy_my_jck"
